$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 13 - Security Value Posting Group: add Danish title
$ws.Range("D13").Value = "Værdipapirværdibogføringsgruppe"

# Row 16 - Security Detailed Information: add Status "Skip?"
$ws.Range("E16").Value = "Skip?"

# Rows 20-22: change "Skip" to "Skip?"
$ws.Range("E20").Value = "Skip?"
$ws.Range("E21").Value = "Skip?"
$ws.Range("E22").Value = "Skip?"

# Row 23: rename ISIN-code -> ISIN, add Danish title "ISIN"
$ws.Range("C23").Value = "ISIN"
$ws.Range("D23").Value = "ISIN"

# Row 24 - Detailed Value Ledger Entry: add Danish title
$ws.Range("D24").Value = "Detaljeret værdipost"

# Row 25 - Detailed Profit Ledger Entry: add Danish title
$ws.Range("D25").Value = "Detaljeret afkastspost"

# Rows 26-29: add Status "Skip?"
$ws.Range("E26").Value = "Skip?"
$ws.Range("E27").Value = "Skip?"
$ws.Range("E28").Value = "Skip?"
$ws.Range("E29").Value = "Skip?"
